# InitProperty.xlsx: unify the conception of DataNode, DataTable, Entity.
#
# The sheet that used to represent a single "Property1" table is renamed to
# the more generic "DataNode", and the in-progress edit position (selection /
# scroll) is moved further down the sheet to around row 50 in column E,
# reflecting where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab: "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Make sure we're looking at the renamed sheet and move the cursor/selection
# to where editing continued (row ~50, column E) while keeping the existing
# frozen header rows (1-8) in place.
$ws.Activate()
$ws.Range("E50").Select()
